# "results of version space"
# - Move the active selection from H1 to F15.
# - Narrow columns B:H (the numeric/result columns) to their new widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes (column letter -> new ColumnWidth in characters).
# Values chosen so the saved OOXML <col .../> width attribute lands on the
# target width from the diff (stored width = ColumnWidth + 5/6, rounded to
# the nearest 1/6 of a character as Excel's column-width grid requires).
$ws.Columns.Item(2).ColumnWidth = 9.66666666666667   # B: 11.4897959183673 -> 10.4948979591837
$ws.Columns.Item(3).ColumnWidth = 9.5                # C: 11.2959183673469 -> 10.2959183673469
$ws.Columns.Item(4).ColumnWidth = 9.66666666666667   # D: 11.4897959183673 -> 10.4948979591837
$ws.Columns.Item(5).ColumnWidth = 13                 # E: 15.8010204081633 -> 13.7959183673469
$ws.Columns.Item(6).ColumnWidth = 12                 # F: 14.7959183673469 -> 12.7959183673469
$ws.Columns.Item(7).ColumnWidth = 12.8333333333333   # G: 15.6938775510204 -> 13.6989795918367
$ws.Columns.Item(8).ColumnWidth = 10.6666666666667   # H: 7.79591836734694 -> 11.4897959183673

# Move the selection / active cell to F15.
$ws.Range("F15").Select()
